$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the "Maximum number of attendee" formula cell to "max_number_of_people"
$ws.Range("D6").Value = '${twig:record.max_number_of_people}'
